# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" values for the 9014c13d... file row across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 12:48:38"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-18 12:48:33"
$wsZhCn.Range("K3").Value = "2016-08-18 12:48:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-18 12:48:38"
$wsDeDe.Range("K3").Value = "2016-08-18 12:48:57"
